$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 72: add meeting note (C72) and hours (E72), grow the row to fit the wrapped text ---
$ws.Range("C72").WrapText = $true
$ws.Range("C72").VerticalAlignment = -4160
$ws.Range("C72").Value = "9:00 to 16:00 => Meeting with Justina; Reworking some of the plots to fit appropriate ranges; Reworking imputations by applying conversions; Designing boxplots instead of average line graphs"
$ws.Range("E72").Value = 7
$ws.Rows.Item(72).RowHeight = 51

# --- Rows 78-87: continue the weekly date/day-of-week table, copying format from the prior week ---
$ws.Range("A71").Copy()
$ws.Range("A78").PasteSpecial(-4122)
$ws.Range("A78").Value = 43863
$ws.Range("B71").Copy()
$ws.Range("B78").PasteSpecial(-4122)
$ws.Range("B78").Value = "Sunday"

$ws.Range("A72").Copy()
$ws.Range("A79").PasteSpecial(-4122)
$ws.Range("A79").Value = 43864
$ws.Range("B72").Copy()
$ws.Range("B79").PasteSpecial(-4122)
$ws.Range("B79").Value = "Monday"

$ws.Range("A73").Copy()
$ws.Range("A80").PasteSpecial(-4122)
$ws.Range("A80").Value = 43865
$ws.Range("B73").Copy()
$ws.Range("B80").PasteSpecial(-4122)
$ws.Range("B80").Value = "Tuesday"

$ws.Range("A74").Copy()
$ws.Range("A81").PasteSpecial(-4122)
$ws.Range("A81").Value = 43866
$ws.Range("B74").Copy()
$ws.Range("B81").PasteSpecial(-4122)
$ws.Range("B81").Value = "Wednesday"

$ws.Range("A75").Copy()
$ws.Range("A82").PasteSpecial(-4122)
$ws.Range("A82").Value = 43867
$ws.Range("B75").Copy()
$ws.Range("B82").PasteSpecial(-4122)
$ws.Range("B82").Value = "Thursday"

$ws.Range("A76").Copy()
$ws.Range("A83").PasteSpecial(-4122)
$ws.Range("A83").Value = 43868
$ws.Range("B76").Copy()
$ws.Range("B83").PasteSpecial(-4122)
$ws.Range("B83").Value = "Friday"

$ws.Range("A77").Copy()
$ws.Range("A84").PasteSpecial(-4122)
$ws.Range("A84").Value = 43869
$ws.Range("B77").Copy()
$ws.Range("B84").PasteSpecial(-4122)
$ws.Range("B84").Value = "Saturday"

$ws.Range("A71").Copy()
$ws.Range("A85").PasteSpecial(-4122)
$ws.Range("A85").Value = 43870
$ws.Range("B71").Copy()
$ws.Range("B85").PasteSpecial(-4122)
$ws.Range("B85").Value = "Sunday"

$ws.Range("A72").Copy()
$ws.Range("A86").PasteSpecial(-4122)
$ws.Range("A86").Value = 43871
$ws.Range("B72").Copy()
$ws.Range("B86").PasteSpecial(-4122)
$ws.Range("B86").Value = "Monday"

$ws.Range("A73").Copy()
$ws.Range("A87").PasteSpecial(-4122)
$ws.Range("A87").Value = 43872
$ws.Range("B73").Copy()
$ws.Range("B87").PasteSpecial(-4122)
$ws.Range("B87").Value = "Tuesday"

$excel.CutCopyMode = $false

# --- Row 85: weekly totals (hours worked, pay) ---
$ws.Range("F85").Formula = "=SUM(E72:E84)"
$ws.Range("G85").Formula = "=F85*20"

# --- View state: scroll down to the new rows and move the selection ---
$excel.ActiveWindow.ScrollRow = 55
$ws.Range("C73").Select()
